# Sync attendance_reports, modules_schedules, and assets from main repo - 2025-12-31 17:14:02
#
# The "Recorded By" column (G) lists the users who recorded/edited a session,
# as a comma-separated string. This script re-normalizes the ordering of
# those comma-separated values for the specific value-sets that changed
# upstream, applying the exact same before -> after text substitution to
# every matching cell in column G, no matter which row it appears on.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of exact "before" text -> exact "after" text for the G (Recorded By) column.
$map = @{
    "System, backup@backdoor.com, system" = "system, System, backup@backdoor.com"
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System"
    "System, admin@admin.com"             = "admin@admin.com, System"
    "admin@admin.com, dnasr281@gmail.com" = "dnasr281@gmail.com, admin@admin.com"
}

$lastRow = $ws.UsedRange.Rows.Count

for ($i = 2; $i -le $lastRow; $i++) {
    $cell = $ws.Cells.Item($i, 7)   # Column G = "Recorded By"
    $current = $cell.Text

    if ($map.ContainsKey($current)) {
        $cell.Value = $map[$current]
    }
}
